$d = $word.ActiveDocument

function Replace-Text($old, $new) {
  $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "Appendix 17: SWIFT Referrals Interview: Information Sheet and Consent Form" "Bylaag 17: SWIFT Verwysingsonderhoud: Inligtingsblad en Toestemmingsvorm"
Replace-Text " What will my interview look like and what is expected of me?" " Hoe sal my onderhoud lyk en wat word van my verwag?"
Replace-Text "Why have I been invited to the interview?" "Waarom is ek na die onderhoud genooi?"
Replace-Text "Do I have to agree to be interviewed?" "Moet ek instem om ondervra te word?"
Replace-Text "What happens with my information?" "Wat gebeur met my inligting?"
Replace-Text "To protect your personal information (including your real name, contact details, and any other information that can identify you), we will give you a participant number, and you can choose a name you want us to call you during the interview." "Om jou persoonlike inligting (insluitend jou regte naam, kontakbesonderhede, en enige ander inligting wat jou kan identifiseer) te beskerm, sal ons vir jou 'n deelnemernommer gee, en jy kan 'n naam kies waarmee ons jou tydens die onderhoud kan aanspreek."
Replace-Text "Do I get anything for being interviewed? " "Kry ek enige iets vir deelname aan die onderhoud? "
Replace-Text "As a thank you for taking part in the discussion, we will give you a R30 airtime voucher/data bundle. " "As 'n bedanking vir jou deelname aan die gesprek, sal ons vir jou 'n R30 lugtydkoepon/data-bundel gee. "
Replace-Text "What happens to my information if I agree to be interviewed?" "Wat gebeur met my inligting as ek instem om ondervra te word?"
Replace-Text "Who are some of the study team members?" "Wie is sommige van die spanlede van die studie?"
Replace-Text "Are there any risks in being interviewed?   " "Is daar enige risiko's verbonde aan die onderhoud?   "
Replace-Text "Who pays for the study?" "Wie betaal vir die studie?"
Replace-Text "This study is part of the Global Parenting Initiative, funded by the LEGO Foundation, Oak Foundation, the World Childhood Foundation, The Human Safety Net, and the UK Research and Innovation Global Challenges Research Fund. " "Hierdie studie is deel van die Global Parenting Initiative, gefinansier deur die LEGO Foundation, Oak Foundation, die World Childhood Foundation, The Human Safety Net, en die UK Research and Innovation Global Challenges Research Fund. "
Replace-Text "Data protection" "Databeskerming"
Replace-Text "Who has approved this study?" "Wie het hierdie studie goedgekeur?"
Replace-Text "Who do I contact if I have questions or concerns?" "Wie kan ek kontak as ek vrae of bekommernisse het?"
Replace-Text "If you have any questions or concerns about your rights as a study participant, you can contact the study team at swift@globalparenting.org or on WhatsApp at +27 XX XXX XXXX (messages only)." "As jy enige vrae of bekommernisse het oor jou regte as 'n studie-deelnemer, kan jy die studiespan kontak by swift@globalparenting.org of via WhatsApp by +27 XX XXX XXXX (net boodskappe)."
Replace-Text "If you have more questions or concerns about your rights, you can contact one of the ethics committees listed: " "As jy meer vrae of bekommernisse het oor jou regte, kan jy een van die etiekkomitees hieronder kontak: "
Replace-Text "Name" "Naam"
Replace-Text "Telephone" "Telefoon"
Replace-Text "Email" "E-pos"
Replace-Text "University of Cape Town Centre for Social Science Research " "Universiteit van Kaapstad Sentrum vir Sosiale Wetenskap Navorsing "
Replace-Text "Human Research Ethics Committee" "Etiekkomitee vir Menslike Navorsing"
Replace-Text "Informed Telephonic consent to take part in the study." "Ingeligte Telefoniese toestemming om aan die studie deel te neem."
Replace-Text "Someone from the research team has gone over all the information above and I know what I need to do." "Iemand van die navorsingspan het al die inligting hierbo deurgegaan en ek weet wat ek moet doen."
Replace-Text "I know who can see my information after the interview, how it will be kept safe, and what happens to it after the study." "Ek weet wie my inligting na die onderhoud kan sien, hoe dit veilig gehou sal word, en wat daarmee sal gebeur na die studie."
Replace-Text "I know I can request access to my data, correct any mistakes, ask to delete it, or for it to be transferred somewhere else." "Ek weet ek kan toegang tot my data versoek, enige foute regstel, vra dat dit verwyder word, of vir dit om na 'n ander plek oorgedra te word."
$rsquo = [char]0x2019
$oldWontBeNamed = "I know that I won" + $rsquo + "t be named in any papers or reports from this study."
Replace-Text $oldWontBeNamed "Ek weet dat ek nie in enige artikels of verslae van hierdie studie genoem sal word nie."
Replace-Text "I know who to tell if I have a problem with the study." "Ek weet wie ek moet kontak as ek 'n probleem met die studie het."
Replace-Text "I can be contacted again if more information is needed from me." "Ek kan weer gekontak word as meer inligting van my nodig is."
Replace-Text "I understand the team will keep my contact information safe so they can tell me about the results of the study." "Ek verstaan dat die span my kontakbesonderhede veilig sal hou sodat hulle my kan inlig oor die resultate van die studie."
